# Updated test query and results files.
#
# The "export" sheet (Worksheets.Item(2), tabSelected) gains three new
# columns - wpid, version, species - inserted between the existing
# "url" and "direct count" columns. The former "direct count" / "direct
# list" columns (C, D) shift right to become F, G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert three blank columns at C:E; existing C/D (direct count/direct
# list) shift to F/G automatically, carrying their values with them.
$ws.Columns("C:E").Insert()

# New header row for the inserted columns.
$ws.Range("C1").Value = "wpid"
$ws.Range("D1").Value = "version"
$ws.Range("E1").Value = "species"

# Per-pathway metadata, in the same row order as the already-sorted
# sheet (rows 2-9).
$wpid    = @("WP2292", "WP2309", "WP232", "WP2841", "WP385", "WP493", "WP539", "WP553")
$version = "WikiPathways_20170210"
$species = "Mus musculus"

for ($i = 0; $i -lt $wpid.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $wpid[$i]
    $ws.Cells.Item($r, 4).Value = $version
    $ws.Cells.Item($r, 5).Value = $species
}

# Match the author's observed column widths for the two new text
# columns that needed extra room.
$ws.Columns.Item(4).ColumnWidth = 21.6640625
$ws.Columns.Item(5).ColumnWidth = 12.83203125

# Refresh the sheet's remembered sort state (still descending on
# "direct count", which is now column F instead of C) to match the
# shifted layout.
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("F1:F9"), 0, 2, 0, 1) | Out-Null
$sort.SetRange($ws.Range("A1:G9"))
$sort.Header = 1
$sort.Apply()

# Selection moved off-sheet to D13 in the committed file.
$ws.Range("D13").Select() | Out-Null

# Keep the page oriented portrait (as in the committed file).
$ws.PageSetup.Orientation = 1 | Out-Null
